$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.4
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 8.5
$ws.Range("J2").Value = 1.95
$ws.Range("K2").Value = 2.3
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 9.5
$ws.Range("AC2").Value = 9.5
$ws.Range("AD2").Value = 8.5
$ws.Range("AG2").Value = 17
$ws.Range("AJ2").Value = 101
$ws.Range("AK2").Value = 67
$ws.Range("G3").Value = 1.67
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 2.38
$ws.Range("N3").Value = 7.5
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("X3").Value = 7
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AH3").Value = 26
$ws.Range("AI3").Value = 19
$ws.Range("AN3").Value = 3.5
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 67
$ws.Range("AZ3").Value = 126
$ws.Range("BA3").Value = 151
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.25
$ws.Range("N5").Value = 9.5
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 1.83
$ws.Range("AL5").Value = 51
$ws.Range("AN5").Value = 3.5
$ws.Range("AW5").Value = 7
$ws.Range("BB5").Value = 301
$ws.Range("N13").Value = 15
$ws.Range("G15").Value = 1.7
$ws.Range("H15").Value = 3.25
$ws.Range("I15").Value = 5
$ws.Range("W15").Value = 5.5
$ws.Range("X15").Value = 7
$ws.Range("Z15").Value = 13
$ws.Range("AH15").Value = 23
$ws.Range("AU15").Value = 9.5
$ws.Range("AW15").Value = 6.5
$ws.Range("AX15").Value = 29
$ws.Range("G20").Value = 1.75
$ws.Range("H20").Value = 3.5
$ws.Range("I20").Value = 4.75
$ws.Range("J20").Value = 2.4
$ws.Range("L20").Value = 5.5
$ws.Range("M20").Value = 1.08
$ws.Range("N20").Value = 8
$ws.Range("O20").Value = 1.36
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 2.2
$ws.Range("R20").Value = 1.65
$ws.Range("S20").Value = 1.44
$ws.Range("T20").Value = 2.63
$ws.Range("U20").Value = 2.1
$ws.Range("V20").Value = 1.67
$ws.Range("W20").Value = 6
$ws.Range("X20").Value = 7.5
$ws.Range("Z20").Value = 13
$ws.Range("AA20").Value = 17
$ws.Range("AD20").Value = 6.5
$ws.Range("AE20").Value = 19
$ws.Range("AG20").Value = 11
$ws.Range("AH20").Value = 23
$ws.Range("AI20").Value = 17
$ws.Range("AK20").Value = 41
$ws.Range("AN20").Value = 3.6
$ws.Range("AO20").Value = 9.5
$ws.Range("AQ20").Value = 34
$ws.Range("AT20").Value = 2.63
$ws.Range("AU20").Value = 9
$ws.Range("AW20").Value = 6.5
$ws.Range("AX20").Value = 29
$ws.Range("AZ20").Value = 101
$ws.Range("BA20").Value = 126
$ws.Range("G21").Value = 2.05
$ws.Range("H21").Value = 3.1
$ws.Range("I21").Value = 3.9
$ws.Range("J21").Value = 2.88
$ws.Range("M21").Value = 1.1
$ws.Range("N21").Value = 7
$ws.Range("O21").Value = 1.44
$ws.Range("P21").Value = 2.63
$ws.Range("Q21").Value = 2.4
$ws.Range("R21").Value = 1.53
$ws.Range("S21").Value = 1.53
$ws.Range("T21").Value = 2.38
$ws.Range("Y21").Value = 9.5
$ws.Range("Z21").Value = 19
$ws.Range("AC21").Value = 7
$ws.Range("AG21").Value = 9
$ws.Range("AK21").Value = 34
$ws.Range("AN21").Value = 4
$ws.Range("AO21").Value = 12
$ws.Range("AP21").Value = 26
$ws.Range("AT21").Value = 2.38
$ws.Range("G23").Value = 2.45
$ws.Range("I23").Value = 2.75
$ws.Range("J23").Value = 3.2
$ws.Range("N23").Value = 9.5
$ws.Range("X23").Value = 12
$ws.Range("AO23").Value = 15
$ws.Range("O25").Value = 1.29
$ws.Range("P25").Value = 3.5
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 1.85
$ws.Range("G32").Value = 2.32
$ws.Range("I32").Value = 2.82
$ws.Range("J32").Value = 2.92
$ws.Range("L32").Value = 3.45
$ws.Range("O32").Value = 1.29
$ws.Range("Q32").Value = 1.87
$ws.Range("R32").Value = 1.87
$ws.Range("T32").Value = 2.72
$ws.Range("W32").Value = 8.75
$ws.Range("Y32").Value = 9
$ws.Range("Z32").Value = 25
$ws.Range("AA32").Value = 18
$ws.Range("AB32").Value = 26
$ws.Range("AE32").Value = 13
$ws.Range("AG32").Value = 9.25
$ws.Range("AI32").Value = 10.25
$ws.Range("AJ32").Value = 35
$ws.Range("AK32").Value = 24
$ws.Range("AL32").Value = 30
$ws.Range("AN32").Value = 4.3
$ws.Range("AO32").Value = 12.5
$ws.Range("AP32").Value = 19.5
$ws.Range("AQ32").Value = 50
$ws.Range("AR32").Value = 80
$ws.Range("AT32").Value = 2.72
$ws.Range("AU32").Value = 6.9
$ws.Range("AV32").Value = 60
$ws.Range("AW32").Value = 4.85
$ws.Range("AX32").Value = 16
$ws.Range("AY32").Value = 23
$ws.Range("AZ32").Value = 75
$ws.Range("BA32").Value = 110
$ws.Range("BB32").Value = 300
$ws.Range("G35").Value = 1.78
$ws.Range("H35").Value = 3.25
$ws.Range("I35").Value = 4.5
$ws.Range("J35").Value = 2.42
$ws.Range("M35").Value = 8.199999999999999
$ws.Range("N35").Value = 1.06
$ws.Range("O35").Value = 1.28
$ws.Range("P35").Value = 3.05
$ws.Range("Q35").Value = 1.88
$ws.Range("R35").Value = 1.82
$ws.Range("S35").Value = 1.4
$ws.Range("T35").Value = 2.5
$ws.Range("U35").Value = 1.72
$ws.Range("V35").Value = 1.9
$ws.Range("W35").Value = 6.7
$ws.Range("X35").Value = 8.25
$ws.Range("Z35").Value = 15
$ws.Range("AB35").Value = 26
$ws.Range("AC35").Value = 9.5
$ws.Range("AD35").Value = 6.4
$ws.Range("AG35").Value = 13.5
$ws.Range("AH35").Value = 29
$ws.Range("AI35").Value = 14
$ws.Range("AJ35").Value = 80
$ws.Range("AL35").Value = 40
$ws.Range("AM35").Value = 450
$ws.Range("AN35").Value = 3.6
$ws.Range("AO35").Value = 9.25
$ws.Range("AP35").Value = 18.5
$ws.Range("AR35").Value = 70
$ws.Range("AW35").Value = 6.1
$ws.Range("AY35").Value = 28
$ws.Range("BA35").Value = 150
$ws.Range("BB35").Value = 350
$ws.Range("M37").Value = 1.02
$ws.Range("N37").Value = 7.1
